$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3225.2222
$ws.Range("I34").Value = 3134.625
$ws.Range("J34").Value = 3950
$ws.Range("K34").Value = 3134.625
$ws.Range("L34").Value = 3950
$ws.Range("M34").Value = -2931.625
$ws.Range("N34").Value = -4356

$ws.Range("H36").Value = 3225.2222
$ws.Range("I36").Value = 3134.625
$ws.Range("J36").Value = 3950
$ws.Range("K36").Value = 3134.625
$ws.Range("L36").Value = 3950
$ws.Range("M36").Value = -2419.625
$ws.Range("N36").Value = -5380

$ws.Range("H111").Value = 1410.5
$ws.Range("I111").Value = 1341.4445
$ws.Range("K111").Value = 4024.3335
$ws.Range("M111").Value = -957.3335000000002

$ws.Range("H112").Value = 2902.8708
$ws.Range("I112").Value = 2409.25
$ws.Range("J112").Value = 3074.5652
$ws.Range("K112").Value = 7227.75
$ws.Range("L112").Value = 9223.695599999999
$ws.Range("M112").Value = -6119.75
$ws.Range("N112").Value = -11439.6956

$ws.Range("H135").Value = 7058
$ws.Range("I135").Value = 860.6667
$ws.Range("J135").Value = 9714
$ws.Range("K135").Value = 7746.0003
$ws.Range("L135").Value = 87426
$ws.Range("M135").Value = -5211.0003
$ws.Range("N135").Value = -92496

$ws.Range("H137").Value = 38469916
$ws.Range("I137").Value = 125001360
$ws.Range("J137").Value = 11494.389
$ws.Range("K137").Value = 375004080
$ws.Range("L137").Value = 34483.167
$ws.Range("M137").Value = -375001530
$ws.Range("N137").Value = -39583.167

$ws.Range("H138").Value = 3727.7742
$ws.Range("J138").Value = 3568.1035
$ws.Range("L138").Value = 10704.3105
$ws.Range("N138").Value = -20984.3105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 662.3333
$ws.Range("I2").Value = 597.58826
$ws.Range("J2").Value = 937.5
$ws.Range("K2").Value = 597.58826
$ws.Range("L2").Value = 937.5
$ws.Range("M2").Value = -484.58826
$ws.Range("N2").Value = -1163.5

$ws.Range("H32").Value = 9258.271000000001
$ws.Range("I32").Value = 4099.7856
$ws.Range("J32").Value = 16996
$ws.Range("K32").Value = 4099.7856
$ws.Range("L32").Value = 16996
$ws.Range("M32").Value = -3812.7856
$ws.Range("N32").Value = -17570

$ws.Range("H74").Value = 3095257.8
$ws.Range("I74").Value = 3971295.5
$ws.Range("K74").Value = 3971295.5
$ws.Range("M74").Value = -3970421.5

$ws.Range("H77").Value = 3095257.8
$ws.Range("I77").Value = 3971295.5
$ws.Range("K77").Value = 19856477.5
$ws.Range("M77").Value = -19852109.5

$ws.Range("H116").Value = 662.3333
$ws.Range("I116").Value = 597.58826
$ws.Range("J116").Value = 937.5
$ws.Range("K116").Value = 597.58826
$ws.Range("L116").Value = 937.5
$ws.Range("M116").Value = 1696.41174
$ws.Range("N116").Value = -5525.5

$ws.Range("H122").Value = 2207.7144
$ws.Range("I122").Value = 1790.8
$ws.Range("K122").Value = 5372.4
$ws.Range("M122").Value = -2922.4

$ws.Range("H132").Value = 612043.4399999999
$ws.Range("I132").Value = 783524.3
$ws.Range("K132").Value = 2350572.9
$ws.Range("M132").Value = -2348042.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 662.3333
$ws.Range("I3").Value = 597.58826
$ws.Range("J3").Value = 937.5
$ws.Range("K3").Value = 597.58826
$ws.Range("L3").Value = 937.5
$ws.Range("M3").Value = -483.58826
$ws.Range("N3").Value = -1165.5

$ws.Range("H122").Value = 58500
$ws.Range("J122").Value = 58500
$ws.Range("L122").Value = 58500
$ws.Range("N122").Value = -68300

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5559629
$ws.Range("I31").Value = 7940942.5
$ws.Range("J31").Value = 3231.3333
$ws.Range("K31").Value = 7940942.5
$ws.Range("L31").Value = 3231.3333
$ws.Range("M31").Value = -7940647.5
$ws.Range("N31").Value = -3821.3333

$ws.Range("H34").Value = 5559629
$ws.Range("I34").Value = 7940942.5
$ws.Range("J34").Value = 3231.3333
$ws.Range("K34").Value = 7940942.5
$ws.Range("L34").Value = 3231.3333
$ws.Range("M34").Value = -7940740.5
$ws.Range("N34").Value = -3635.3333

$ws.Range("H99").Value = 35419.445
$ws.Range("I99").Value = 45799.8
$ws.Range("J99").Value = 22444
$ws.Range("K99").Value = 45799.8
$ws.Range("L99").Value = 22444
$ws.Range("M99").Value = -44301.8
$ws.Range("N99").Value = -25440

$ws.Range("H103").Value = 8940.714
$ws.Range("I103").Value = 8940.714
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 8940.714
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -7768.714
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 21535.666
$ws.Range("I122").Value = 4161.4
$ws.Range("J122").Value = 43253.5
$ws.Range("K122").Value = 12484.2
$ws.Range("L122").Value = 129760.5
$ws.Range("M122").Value = -10034.2
$ws.Range("N122").Value = -134660.5

$ws.Range("H126").Value = 35419.445
$ws.Range("I126").Value = 45799.8
$ws.Range("J126").Value = 22444
$ws.Range("K126").Value = 137399.4
$ws.Range("L126").Value = 67332
$ws.Range("M126").Value = -134929.4
$ws.Range("N126").Value = -72272

$ws.Range("H141").Value = 149437.78
$ws.Range("J141").Value = 170659.4
$ws.Range("L141").Value = 170659.4
$ws.Range("N141").Value = -181019.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 312619.88
$ws.Range("I14").Value = 312619.88
$ws.Range("K14").Value = 937859.64
$ws.Range("M14").Value = -937686.64

$ws.Range("H117").Value = 1291.8889
$ws.Range("I117").Value = 519.5714
$ws.Range("J117").Value = 3995
$ws.Range("K117").Value = 1558.7142
$ws.Range("L117").Value = 11985
$ws.Range("M117").Value = 1883.2858
$ws.Range("N117").Value = -18869

$ws.Range("H136").Value = 5914.8887
$ws.Range("I136").Value = 2536.6667
$ws.Range("J136").Value = 8617.467000000001
$ws.Range("K136").Value = 7610.000100000001
$ws.Range("L136").Value = 25852.401
$ws.Range("M136").Value = -2510.000100000001
$ws.Range("N136").Value = -36052.401

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1639
$ws.Range("I97").Value = 1606.7858
$ws.Range("J97").Value = 1789.3334
$ws.Range("K97").Value = 1606.7858
$ws.Range("L97").Value = 1789.3334
$ws.Range("M97").Value = -1110.7858
$ws.Range("N97").Value = -2781.3334

$ws.Range("H122").Value = 63172.95
$ws.Range("I122").Value = 82413.36
$ws.Range("J122").Value = 9299.799999999999
$ws.Range("K122").Value = 247240.08
$ws.Range("L122").Value = 27899.4
$ws.Range("M122").Value = -244790.08
$ws.Range("N122").Value = -32799.39999999999

$ws.Range("H126").Value = 2671.1428
$ws.Range("I126").Value = 2459.6
$ws.Range("K126").Value = 7378.799999999999
$ws.Range("M126").Value = -4908.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11007.111
$ws.Range("I7").Value = 15471.167
$ws.Range("J7").Value = 2079
$ws.Range("K7").Value = 15471.167
$ws.Range("L7").Value = 2079
$ws.Range("M7").Value = -15359.167
$ws.Range("N7").Value = -2303

$ws.Range("H40").Value = 4051.7368
$ws.Range("I40").Value = 3686.875
$ws.Range("K40").Value = 3686.875
$ws.Range("M40").Value = -3550.875

$ws.Range("H46").Value = 3704.8823
$ws.Range("I46").Value = 490.75
$ws.Range("K46").Value = 490.75
$ws.Range("M46").Value = -302.75

$ws.Range("H55").Value = 1415.5652
$ws.Range("I55").Value = 1357.909
$ws.Range("K55").Value = 1357.909
$ws.Range("M55").Value = -1184.909

$ws.Range("H74").Value = 49938.2
$ws.Range("I74").Value = 49938.2
$ws.Range("K74").Value = 49938.2
$ws.Range("M74").Value = -48940.2

$ws.Range("H77").Value = 49938.2
$ws.Range("I77").Value = 49938.2
$ws.Range("K77").Value = 149814.6
$ws.Range("M77").Value = -144822.6

$ws.Range("H122").Value = 4649
$ws.Range("I122").Value = 2976
$ws.Range("K122").Value = 8928
$ws.Range("M122").Value = -6478

$ws.Range("H126").Value = 11007.111
$ws.Range("I126").Value = 15471.167
$ws.Range("J126").Value = 2079
$ws.Range("K126").Value = 46413.501
$ws.Range("L126").Value = 6237
$ws.Range("M126").Value = -43943.501
$ws.Range("N126").Value = -11177

$ws.Range("H132").Value = 6876628.5
$ws.Range("I132").Value = 14608505
$ws.Range("K132").Value = 43825515
$ws.Range("M132").Value = -43822985

$ws.Range("H133").Value = 77126.14
$ws.Range("J133").Value = 69496.5
$ws.Range("L133").Value = 69496.5
$ws.Range("N133").Value = -74556.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954

$ws.Range("H122").Value = 70818
$ws.Range("I122").Value = 5545.6924
$ws.Range("J122").Value = 282953
$ws.Range("K122").Value = 16637.0772
$ws.Range("L122").Value = 848859
$ws.Range("M122").Value = -14187.0772
$ws.Range("N122").Value = -853759

$ws.Range("H132").Value = 6947909
$ws.Range("I132").Value = 7939779.5
$ws.Range("J132").Value = 4816.3335
$ws.Range("K132").Value = 23819338.5
$ws.Range("L132").Value = 14449.0005
$ws.Range("M132").Value = -23816808.5
$ws.Range("N132").Value = -19509.0005
